$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "Eric Carman" -> "Eric Cartman" (row 3)
$ws.Cells.Item(3, 1).Value = "Eric Cartman"

# Insert a new row 4 for Peggy Hill's data (shifts old rows 4-5 down to 5-6)
$ws.Rows("4:4").Insert()
$ws.Cells.Item(4, 1).Value = "Peggy Hill"
$ws.Cells.Item(4, 2).Value = 38
$ws.Cells.Item(4, 3).Value = 66
$ws.Cells.Item(4, 4).Value = 125
$ws.Cells.Item(4, 5).Value = "Hank Hill"

# Add a new row 7 for Marge Simpson's data
$ws.Cells.Item(7, 1).Value = "Marge Simpson"
$ws.Cells.Item(7, 2).Value = 43
$ws.Cells.Item(7, 3).Value = 78
$ws.Cells.Item(7, 4).Value = 135
$ws.Cells.Item(7, 5).Value = "Homer J. Simpson"

# Update the selection to match the final state
[void]$ws.Range("C10").Select()
